# Update "想去人数" (F column) figures on both the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 3..6 hold F values that need bumping
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 292
$wsExhibit.Range("F4").Value = 1285
$wsExhibit.Range("F5").Value = 81
$wsExhibit.Range("F6").Value = 60

# Sheet "全部类型" - rows 4..7 hold the same events, shifted by one row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 292
$wsAll.Range("F5").Value = 1285
$wsAll.Range("F6").Value = 81
$wsAll.Range("F7").Value = 60
